$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Attributs" paragraph: extend the list of known magic schools.
#    (Assigning Range.Text directly -- instead of using Find.Execute's
#    replacement argument -- keeps straight apostrophes as-is; the Execute
#    replacement path "smart-quotes" them, which we don't want here.)
# ---------------------------------------------------------------------------
$rAttr = $d.Content
$rAttr.Find.Execute("Magie Aquatique — magie de l'eau, Corps artificiels: créature artificielle, nul besoin pour elle de respirer")
if ($rAttr.Find.Found) {
    $rAttr.Text = "Magie Aquatique — magie de l'eau, Magie Céleste — magie du ciel, Magie Démoniaque — magie liée aux ténèbres, Magie Divine — magie liée aux divinités, Magie Ignis — magie du feu, Magie Naturelle — magie de la nature, Magie Neutre — magie neutre, Magie Terrestre: magie de la terre"
}

# ---------------------------------------------------------------------------
# 2) "Inventaires" list: drop the Scramasax / Épée courte / Épée longue /
#    Glaive entries (whole paragraphs) and also drop the "Couteau..." run,
#    while keeping the trailing <br/> + "Sortilèges : " run that shares the
#    last paragraph with "Couteau...".
# ---------------------------------------------------------------------------

# 2a) remove the four full paragraphs "Scramasax" .. "Glaive" (including
#     their paragraph marks) in one shot.
$rStart = $d.Content
$rStart.Find.Execute("Scramasax 0.8kg 60 cm 0 Tranchant 1d4+1")
$startPos = $rStart.Start

$rEnd = $d.Content
$rEnd.Find.Execute("Glaive 0.8kg 60 cm 0 Tranchant 1d4+1")
$endPos = $rEnd.End + 1

$d.Range($startPos, $endPos).Delete()

# 2b) remove just the "Couteau..." text run (leaving the <br/> and the
#     "Sortilèges : " run alone in their paragraph).
$rCouteau = $d.Content
$rCouteau.Find.Execute("Couteau 0.8kg 60 cm 0 Tranchant 1d4+1")
$rCouteau.Delete()

# ---------------------------------------------------------------------------
# 3) Drop the "Manipulation de l'eau" run entirely, leaving its paragraph
#    empty.
# ---------------------------------------------------------------------------
$rSpell = $d.Content
$rSpell.Find.Execute("Manipulation de l'eau")
$rSpell.Delete()
